# Apply the "Error Calculations and Plots" edits to the missing_data sheet.
#
# Summary of what changed between before/after:
#  - A handful of individual cells in rows 2-25 either became newly "missing"
#    (cleared, matching the workbook's existing convention for missing data)
#    or were newly filled in with a numeric value.
#  - The "RM 232" row (row 26) and the "SC 92" row (originally row 28) were
#    both removed entirely, shifting every following row up.
#  - After that shift, a few cells in the shifted rows (now SC 101, SC 105,
#    SC 119, SC 120, SC 193) also had values filled in or cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell edits in rows 2-25 (row positions unaffected by the later deletes) ---

# E2: -7.2 -> missing
$ws.Range("E2").ClearContents()

# E5: missing -> -5
$ws.Range("E5").Value = -5

# C6: missing -> 15.1, E6: missing -> -5.7
$ws.Range("C6").Value = 15.1
$ws.Range("E6").Value = -5.7

# C8: 15.5 -> missing
$ws.Range("C8").ClearContents()

# E9: -6.8 -> missing
$ws.Range("E9").ClearContents()

# E10: -6.1 -> missing
$ws.Range("E10").ClearContents()

# C12: missing -> 12.5
$ws.Range("C12").Value = 12.5

# C14: 14.4 -> missing
$ws.Range("C14").ClearContents()

# C17: missing -> 11.2
$ws.Range("C17").Value = 11.2

# C18: missing -> 11.5
$ws.Range("C18").Value = 11.5

# C19: 13.2 -> missing
$ws.Range("C19").ClearContents()

# C20: 12.5 -> missing
$ws.Range("C20").ClearContents()

# C23: missing -> 12.2
$ws.Range("C23").Value = 12.2

# E24: missing -> -8.1
$ws.Range("E24").Value = -8.1

# --- Remove the "RM 232" row and the "SC 92" row entirely ---
# Row 26 is "RM 232"; deleting it shifts "SC 92" up to row 27.
$ws.Rows(26).Delete()
# Row 27 is now "SC 92"; delete it too, shifting the rest up again.
$ws.Rows(27).Delete()

# --- Cell edits in the shifted rows (now SC 101, SC 105, SC 119, SC 120, SC 193) ---

# SC 101 is now row 27
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()

# SC 105 is now row 28
$ws.Range("B28").ClearContents()
$ws.Range("E28").ClearContents()

# SC 119 is now row 29
$ws.Range("B29").ClearContents()

# SC 120 is now row 30
$ws.Range("B30").Value = -19.7
$ws.Range("E30").Value = -5.7

# SC 193 is now row 32
$ws.Range("B32").ClearContents()
